# Apply the 2024-04-03 cryptos data refresh to Sheet1 (rows 2-51).
# Prices/volumes were re-scraped from coinranking.com; a handful of rows
# also shifted coin/link because of new ranking order (new row inserted
# for ImmutableX, Chainlink/WrappedBTC swapped order, etc).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = "66.192.16"
$ws.Range("E2").Value = "  -0.58%  "

# Row 3: Ethereum
$ws.Range("D3").Value = "3.309.91"
$ws.Range("E3").Value = "  -1.68%  "

# Row 4: TetherUSD
$ws.Range("E4").Value = "  +0.15%  "

# Row 5: Solana
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "189.92"
$ws.Range("E5").Value = "  +3.55%  "

# Row 6: BNB
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "561.05"
$ws.Range("E6").Value = "  -0.04%  "

# Row 7: USDC
$ws.Range("E7").Value = "  -0.09%  "

# Row 8: XRP
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.588"
$ws.Range("E8").Value = "  -2.12%  "

# Row 9: LidoStakedEther
$ws.Range("D9").Value = "3.299.11"
$ws.Range("E9").Value = "  -1.74%  "

# Row 10: Dogecoin
$ws.Range("E10").Value = "  -1.54%  "

# Row 11: Cardano
$ws.Range("E11").Value = "  -1.16%  "

# Row 12: Avalanche
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "47.77"
$ws.Range("E12").Value = "  -0.56%  "

# Row 13: ShibaInu
$ws.Range("E13").Value = "  +0.54%  "

# Row 14: Polkadot
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.67"
$ws.Range("E14").Value = "  -0.84%  "

# Row 15: WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "3.832.09"
$ws.Range("E15").Value = "  -1.58%  "

# Row 16: BitcoinCash
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "604.92"
$ws.Range("E16").Value = "  -0.58%  "

# Row 17: Chainlink
$ws.Range("B17").Value = "Chainlink"
$ws.Range("C17").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.11"
$ws.Range("E17").Value = "  -0.49%  "

# Row 18: WrappedBTC
$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").Value = "66.180.30"
$ws.Range("E18").Value = "  -0.31%  "

# Row 19: TRON
$ws.Range("E19").Value = "  +0.27%  "

# Row 20: WrappedEther
$ws.Range("D20").Value = "3.303.90"
$ws.Range("E20").Value = "  -1.76%  "

# Row 21: Uniswap
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.11"
$ws.Range("E21").Value = "  -3.68%  "

# Row 22: Polygon
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.914"
$ws.Range("E22").Value = "  -0.26%  "

# Row 23: InternetComputer(DFINITY)
$ws.Range("E23").Value = "  +9.12%  "

# Row 24: Toncoin
$ws.Range("E24").Value = "  -0.90%  "

# Row 25: Litecoin
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "101.12"
$ws.Range("E25").Value = "  +1.23%  "

# Row 26: PancakeSwap
$ws.Range("E26").Value = "  -2.36%  "

# Row 27: ImmutableX
$ws.Range("B27").Value = "ImmutableX"
$ws.Range("C27").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.75"
$ws.Range("E27").Value = "  +0.95%  "

# Row 28: RenderToken
$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.72"
$ws.Range("E28").Value = "  +2.99%  "

# Row 29: Filecoin
$ws.Range("B29").Value = "Filecoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.68"
$ws.Range("E29").Value = "  -1.61%  "

# Row 30: EthereumClassic
$ws.Range("B30").Value = "EthereumClassic"
$ws.Range("C30").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "30.34"
$ws.Range("E30").Value = "  -1.37%  "

# Row 31: NEARProtocol
$ws.Range("B31").Value = "NEARProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.78"
$ws.Range("E31").Value = "  +7.13%  "

# Row 32: dogwifhat
$ws.Range("B32").Value = "dogwifhat"
$ws.Range("C32").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.10"
$ws.Range("E32").Value = "  +6.27%  "

# Row 33: Cosmos
$ws.Range("B33").Value = "Cosmos"
$ws.Range("C33").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "11.13"
$ws.Range("E33").Value = "  -0.08%  "

# Row 34: Bittensor
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "559.77"
$ws.Range("E34").Value = "  +1.53%  "

# Row 35: Hedera
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.106"
$ws.Range("E35").Value = "  -0.24%  "

# Row 36: Dai
$ws.Range("B36").Value = "Dai"
$ws.Range("C36").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("E36").Value = "  +0.07%  "

# Row 37: OKB
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "57.18"
$ws.Range("E37").Value = "  -1.72%  "

# Row 38: Maker
$ws.Range("B38").Value = "Maker"
$ws.Range("C38").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D38").Value = "3.708.17"
$ws.Range("E38").Value = "  -3.82%  "

# Row 39: PEPE
$ws.Range("B39").Value = "PEPE"
$ws.Range("C39").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D39").Value = "0.0₃0729"
$ws.Range("E39").Value = "  +0.60%  "

# Row 40: InjectiveProtocol
$ws.Range("B40").Value = "InjectiveProtocol"
$ws.Range("C40").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "34.15"
$ws.Range("E40").Value = "  +5.47%  "

# Row 41: Stacks
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.31"
$ws.Range("E41").Value = "  -3.23%  "

# Row 42: Kaspa
$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.130"
$ws.Range("E42").Value = "  +1.72%  "

# Row 43: CoreDAO
$ws.Range("B43").Value = "CoreDAO"
$ws.Range("C43").Value = "https://coinranking.com/coin/HFvoXUQh4+coredao-core"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.42"
$ws.Range("E43").Value = "  -1.14%  "

# Row 44: Fetch.AI
$ws.Range("E44").Value = "  +1.28%  "

# Row 45: TheGraph
$ws.Range("B45").Value = "TheGraph"
$ws.Range("C45").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.342"
$ws.Range("E45").Value = "  -1.97%  "

# Row 46: VeChain
$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0423"
$ws.Range("E46").Value = "  +1.50%  "

# Row 47: ApeXProtocol
$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.25"
$ws.Range("E47").Value = "  +3.23%  "

# Row 48: Stellar
$ws.Range("B48").Value = "Stellar"
$ws.Range("C48").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.130"
$ws.Range("E48").Value = "  -0.64%  "

# Row 49: ThetaToken
$ws.Range("B49").Value = "ThetaToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.60"
$ws.Range("E49").Value = "  -2.33%  "

# Row 50: FirstDigitalUSD
$ws.Range("B50").Value = "FirstDigitalUSD"
$ws.Range("C50").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.999"
$ws.Range("E50").Value = "  +0.16%  "

# Row 51: Mantle
$ws.Range("B51").Value = "Mantle"
$ws.Range("C51").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.28"
$ws.Range("E51").Value = "  +1.50%  "
